# MedicationCatalog.pptx - "updated figure for catalog of medications"
#
# 1) Refresh the cached "datetimeFigureOut" footer date (07/05/2021 -> 17/06/2021)
#    on the slide master, every slide layout, and the notes master.
# 2) Split the "Precaution profile of ClinicalUseIssue" caption into two
#    stacked lines: "Warning " / "profile of ClinicalUseIssue".
# 3) Add a new "ingredient.item" label textbox to the diagram.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh
# ---------------------------------------------------------------------
$newDate = "17/06/2021"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.Name -like "*date*" -or $sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 2) Split "Precaution profile of ClinicalUseIssue" into two lines
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 143") {
        $sh.TextFrame.TextRange.Text = "Warning `rprofile of ClinicalUseIssue"
    }
}

# ---------------------------------------------------------------------
# 3) New "ingredient.item" label textbox
# ---------------------------------------------------------------------
$left = 3839844 / 12700.0
$top = 4032675 / 12700.0
$width = 950237 / 12700.0
$height = 230832 / 12700.0

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "ZoneTexte 114"
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "ingredient.item"
$tb.TextFrame.TextRange.Font.Size = 9
$tb.TextFrame.TextRange.Font.Color.RGB = 12611584
$tb.Width = $width
